# Auto-generated script applying scheduled-runner market/profit data refresh
# to the Malboro_Profits leve-profit tables across all job sheets.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6219.3335
$ws.Range("I51").Value = 11250
$ws.Range("J51").Value = 4782
$ws.Range("K51").Value = 11250
$ws.Range("L51").Value = 4782
$ws.Range("M51").Value = -10766
$ws.Range("N51").Value = -5750
$ws.Range("H53").Value = 6733.9473
$ws.Range("I53").Value = 982.44446
$ws.Range("J53").Value = 11910.3
$ws.Range("K53").Value = 982.44446
$ws.Range("L53").Value = 11910.3
$ws.Range("M53").Value = -345.44446
$ws.Range("N53").Value = -13184.3
$ws.Range("H86").Value = 24630.6
$ws.Range("J86").Value = 28941
$ws.Range("L86").Value = 28941
$ws.Range("N86").Value = -31187
$ws.Range("H89").Value = 24630.6
$ws.Range("J89").Value = 28941
$ws.Range("L89").Value = 144705
$ws.Range("N89").Value = -155937
$ws.Range("H111").Value = 529.1053000000001
$ws.Range("I111").Value = 297.73334
$ws.Range("K111").Value = 893.20002
$ws.Range("M111").Value = 2173.79998
$ws.Range("H113").Value = 13355.4375
$ws.Range("I113").Value = 19600.6
$ws.Range("J113").Value = 2946.8333
$ws.Range("K113").Value = 19600.6
$ws.Range("L113").Value = 2946.8333
$ws.Range("M113").Value = -16346.6
$ws.Range("N113").Value = -9454.8333
$ws.Range("H132").Value = 12165.5
$ws.Range("I132").Value = 4999.6665
$ws.Range("J132").Value = 19331.334
$ws.Range("K132").Value = 14998.9995
$ws.Range("L132").Value = 57994.00199999999
$ws.Range("M132").Value = -12468.9995
$ws.Range("N132").Value = -63054.00199999999
$ws.Range("H138").Value = 3647.7932
$ws.Range("I138").Value = 3279.9524
$ws.Range("J138").Value = 3856.5676
$ws.Range("K138").Value = 9839.8572
$ws.Range("L138").Value = 11569.7028
$ws.Range("M138").Value = -4699.8572
$ws.Range("N138").Value = -21849.7028
$ws.Range("H141").Value = 3133.4614
$ws.Range("J141").Value = 4129.1665
$ws.Range("L141").Value = 12387.4995
$ws.Range("N141").Value = -22747.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H102").Value = 6013.1763
$ws.Range("I102").Value = 6013.1763
$ws.Range("K102").Value = 6013.1763
$ws.Range("M102").Value = -4391.1763
$ws.Range("H132").Value = 5764994
$ws.Range("I132").Value = 7823.1875
$ws.Range("J132").Value = 18924242
$ws.Range("K132").Value = 23469.5625
$ws.Range("L132").Value = 56772726
$ws.Range("M132").Value = -20939.5625
$ws.Range("N132").Value = -56777786

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 15672.656
$ws.Range("I134").Value = 10232.036
$ws.Range("K134").Value = 30696.108
$ws.Range("M134").Value = -28161.108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 48994.98
$ws.Range("I31").Value = 68255.266
$ws.Range("J31").Value = 21480.285
$ws.Range("K31").Value = 68255.266
$ws.Range("L31").Value = 21480.285
$ws.Range("M31").Value = -67960.266
$ws.Range("N31").Value = -22070.285
$ws.Range("H34").Value = 48994.98
$ws.Range("I34").Value = 68255.266
$ws.Range("J34").Value = 21480.285
$ws.Range("K34").Value = 68255.266
$ws.Range("L34").Value = 21480.285
$ws.Range("M34").Value = -68053.266
$ws.Range("N34").Value = -21884.285
$ws.Range("H58").Value = 18996.291
$ws.Range("J58").Value = 23832.75
$ws.Range("L58").Value = 23832.75
$ws.Range("N58").Value = -24238.75
$ws.Range("H99").Value = 7064.5
$ws.Range("I99").Value = 6492.6665
$ws.Range("J99").Value = 7750.7
$ws.Range("K99").Value = 6492.6665
$ws.Range("L99").Value = 7750.7
$ws.Range("M99").Value = -4994.6665
$ws.Range("N99").Value = -10746.7
$ws.Range("H107").Value = 1835708.1
$ws.Range("I107").Value = 2202649.8
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2202649.8
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -2200729.8
$ws.Range("N107").Value = -4840
$ws.Range("H126").Value = 7064.5
$ws.Range("I126").Value = 6492.6665
$ws.Range("J126").Value = 7750.7
$ws.Range("K126").Value = 19477.9995
$ws.Range("L126").Value = 23252.1
$ws.Range("M126").Value = -17007.9995
$ws.Range("N126").Value = -28192.1
$ws.Range("H132").Value = 35503700
$ws.Range("I132").Value = 3688.625
$ws.Range("K132").Value = 11065.875
$ws.Range("M132").Value = -8535.875
$ws.Range("H136").Value = 18996.291
$ws.Range("J136").Value = 23832.75
$ws.Range("L136").Value = 71498.25
$ws.Range("N136").Value = -76598.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 61.4
$ws.Range("I2").Value = 60.833332
$ws.Range("J2").Value = 62.25
$ws.Range("K2").Value = 364.999992
$ws.Range("L2").Value = 373.5
$ws.Range("M2").Value = -251.999992
$ws.Range("N2").Value = -599.5
$ws.Range("H5").Value = 464.7619
$ws.Range("I5").Value = 463
$ws.Range("J5").Value = 466.36365
$ws.Range("K5").Value = 1389
$ws.Range("L5").Value = 1399.09095
$ws.Range("M5").Value = -1277
$ws.Range("N5").Value = -1623.09095
$ws.Range("H107").Value = 6109.4736
$ws.Range("J107").Value = 10089.182
$ws.Range("L107").Value = 30267.546
$ws.Range("N107").Value = -34107.546
$ws.Range("H122").Value = 11956958
$ws.Range("J122").Value = 3542411.2
$ws.Range("L122").Value = 31881700.8
$ws.Range("N122").Value = -31886600.8
$ws.Range("H135").Value = 464.7619
$ws.Range("I135").Value = 463
$ws.Range("J135").Value = 466.36365
$ws.Range("K135").Value = 4167
$ws.Range("L135").Value = 4197.27285
$ws.Range("M135").Value = -1632
$ws.Range("N135").Value = -9267.272850000001
$ws.Range("H139").Value = 1227.1875
$ws.Range("I139").Value = 1227.1875
$ws.Range("K139").Value = 3681.5625
$ws.Range("M139").Value = 1458.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 619.6875
$ws.Range("I107").Value = 528.1667
$ws.Range("J107").Value = 674.6
$ws.Range("K107").Value = 528.1667
$ws.Range("L107").Value = 674.6
$ws.Range("M107").Value = 1391.8333
$ws.Range("N107").Value = -4514.6
$ws.Range("H108").Value = 79000
$ws.Range("J108").Value = 79000
$ws.Range("L108").Value = 79000
$ws.Range("M108").Value = -86680
$ws.Range("H113").Value = 2111.8667
$ws.Range("I113").Value = 2165.3
$ws.Range("J113").Value = 2005
$ws.Range("K113").Value = 2165.3
$ws.Range("L113").Value = 2005
$ws.Range("M113").Value = 4.699999999999818
$ws.Range("N113").Value = -6345
$ws.Range("H126").Value = 10319.608
$ws.Range("I126").Value = 8840.362999999999
$ws.Range("K126").Value = 26521.089
$ws.Range("M126").Value = -24051.089
$ws.Range("H132").Value = 441597.03
$ws.Range("I132").Value = 4604.478
$ws.Range("K132").Value = 13813.434
$ws.Range("M132").Value = -11283.434

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2962.5
$ws.Range("I16").Value = 2962.5
$ws.Range("K16").Value = 2962.5
$ws.Range("M16").Value = -2792.5
$ws.Range("H40").Value = 4033.743
$ws.Range("I40").Value = 2827.36
$ws.Range("J40").Value = 7049.7
$ws.Range("K40").Value = 2827.36
$ws.Range("L40").Value = 7049.7
$ws.Range("M40").Value = -2691.36
$ws.Range("N40").Value = -7321.7
$ws.Range("H132").Value = 1558334.4
$ws.Range("I132").Value = 4220.6816
$ws.Range("K132").Value = 12662.0448
$ws.Range("M132").Value = -10132.0448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1569.8334
$ws.Range("I100").Value = 639
$ws.Range("K100").Value = 1278
$ws.Range("M100").Value = -737
$ws.Range("H122").Value = 4431.5947
$ws.Range("I122").Value = 3154.0908
$ws.Range("K122").Value = 9462.2724
$ws.Range("M122").Value = -7012.2724
$ws.Range("H136").Value = 712898.25
$ws.Range("I136").Value = 8001
$ws.Range("K136").Value = 24003
$ws.Range("M136").Value = -21453
